$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 486.6
$ws.Range("I28").Value = 456.08694
$ws.Range("K28").Value = 456.08694
$ws.Range("M28").Value = 28.91305999999997
$ws.Range("H38").Value = 35.25
$ws.Range("I38").Value = 35.25
$ws.Range("K38").Value = 105.75
$ws.Range("M38").Value = 266.25
$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 2500
$ws.Range("K74").Value = 2500
$ws.Range("M74").Value = -1564
$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 2500
$ws.Range("K77").Value = 12500
$ws.Range("M77").Value = -7820
$ws.Range("H98").Value = 977.1177
$ws.Range("I98").Value = 707.4666999999999
$ws.Range("K98").Value = 707.4666999999999
$ws.Range("M98").Value = 790.5333000000001
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H107").Value = 680.7778
$ws.Range("I107").Value = 666.125
$ws.Range("K107").Value = 666.125
$ws.Range("M107").Value = 1253.875
$ws.Range("H122").Value = 977.1177
$ws.Range("I122").Value = 707.4666999999999
$ws.Range("K122").Value = 2122.4001
$ws.Range("M122").Value = 327.5999000000002
$ws.Range("H137").Value = 2036.3334
$ws.Range("J137").Value = 2108.5
$ws.Range("L137").Value = 6325.5
$ws.Range("N137").Value = -11425.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 143.5
$ws.Range("I22").Value = 143.5
$ws.Range("K22").Value = 143.5
$ws.Range("M22").Value = 29.5
$ws.Range("H62").Value = 50001
$ws.Range("J62").Value = 50001
$ws.Range("L62").Value = 50001
$ws.Range("N62").Value = -51373
$ws.Range("H65").Value = 50001
$ws.Range("J65").Value = 50001
$ws.Range("L65").Value = 150003
$ws.Range("N65").Value = -156867
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1420
$ws.Range("N107").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1227.2727
$ws.Range("I7").Value = 849.2857
$ws.Range("K7").Value = 849.2857
$ws.Range("M7").Value = -736.2857
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H43").Value = 34499.5
$ws.Range("J43").Value = 34499.5
$ws.Range("L43").Value = 34499.5
$ws.Range("N43").Value = -34867.5
$ws.Range("H99").Value = 4591.5557
$ws.Range("J99").Value = 5999.5
$ws.Range("L99").Value = 5999.5
$ws.Range("N99").Value = -8995.5
$ws.Range("H101").Value = 34499.5
$ws.Range("J101").Value = 34499.5
$ws.Range("L101").Value = 34499.5
$ws.Range("N101").Value = -40989.5
$ws.Range("H126").Value = 4591.5557
$ws.Range("J126").Value = 5999.5
$ws.Range("L126").Value = 17998.5
$ws.Range("N126").Value = -22938.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 250.73334
$ws.Range("I2").Value = 81.666664
$ws.Range("J2").Value = 363.44446
$ws.Range("K2").Value = 81.666664
$ws.Range("L2").Value = 363.44446
$ws.Range("M2").Value = 31.333336
$ws.Range("N2").Value = -589.4444599999999
$ws.Range("H132").Value = 4024.8572
$ws.Range("I132").Value = 3412.5
$ws.Range("K132").Value = 10237.5
$ws.Range("M132").Value = -7707.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2466.1667
$ws.Range("I7").Value = 3732.3333
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 3732.3333
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = -3620.3333
$ws.Range("N7").Value = -1424
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -3812
$ws.Range("N46").Value = -4376
$ws.Range("H61").Value = 600
$ws.Range("I61").Value = 600
$ws.Range("K61").Value = 600
$ws.Range("M61").Value = -398
$ws.Range("H113").Value = 600
$ws.Range("I113").Value = 600
$ws.Range("K113").Value = 600
$ws.Range("M113").Value = 1570
$ws.Range("H126").Value = 2466.1667
$ws.Range("I126").Value = 3732.3333
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 11196.9999
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -8726.999899999999
$ws.Range("N126").Value = -8540

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -384
$ws.Range("H75").Value = 17129.166
$ws.Range("J75").Value = 16155
$ws.Range("L75").Value = 16155
$ws.Range("N75").Value = -18027
$ws.Range("H78").Value = 17129.166
$ws.Range("J78").Value = 16155
$ws.Range("L78").Value = 48465
$ws.Range("N78").Value = -57825
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H113").Value = 610.25
$ws.Range("I113").Value = 662.3
$ws.Range("K113").Value = 1986.9
$ws.Range("M113").Value = 183.1000000000001
$ws.Range("H122").Value = 1683.1666
$ws.Range("I122").Value = 1739.8
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 5219.4
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -2769.4
$ws.Range("N122").Value = -9100
$ws.Range("H130").Value = 49887.332
$ws.Range("J130").Value = 49887.332
$ws.Range("L130").Value = 49887.332
$ws.Range("N130").Value = -59927.332
$ws.Range("H136").Value = 1120.6471
$ws.Range("I136").Value = 1118
$ws.Range("K136").Value = 3354
$ws.Range("M136").Value = -804
